# Update the confusion-matrix values with corrected translations
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "17 (0.8095)"
$ws.Range("D3").Value = "4 (0.1905)"
$ws.Range("C4").Value = "1 (0.0714)"
$ws.Range("D4").Value = "13 (0.9286)"

# Update the active selection to match the saved view state
$null = $ws.Range("H4").Select()
